# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# --- Part 1: add a new "2022-Q1" sheet. Clone it from "2021-Q4" which
# already carries the identical headers/formatting, then overwrite the few
# data values that differ for the new quarter.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy([System.Reflection.Missing]::Value, $q4)
$newSheet = $wb.Worksheets.Item($q4.Index + 1)
$newSheet.Name = "2022-Q1"

# D2:G2 are stored as text in the source data (keeps exact formatting, e.g.
# the trailing zero in "0.2190"), so format the cells as text before typing
# the values in, same as a user would via Format Cells > Text.
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "5.24"
$newSheet.Range("E2").Value = "46.02"
$newSheet.Range("F2").Value = "4.18"
$newSheet.Range("G2").Value = "0.2190"
$newSheet.Range("H2").Value = 3

# --- Part 2: update the "总计" summary sheet - insert a new top data row for
# 2022-Q1 and shift the existing rows down.
$zj = $wb.Worksheets.Item("总计")
$zj.Rows("2:2").Insert()
$zj.Range("B2:D2").ClearFormats()

# Re-apply the index-column formatting (bold + border) to the new A2 cell by
# copying it from a cell that already carries it.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0.22

# Renumber the index column for the rows that shifted down.
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3
$zj.Range("A6").Value = 4
$zj.Range("A7").Value = 5

# Restore the originally active sheet/tab.
$wb.Worksheets.Item(1).Activate()
